# Apply the edits captured in the commit diff:
#  - rename the single worksheet from "Sheet1" to "Excel"
#  - move the active selection on that sheet from C7 to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (was "Sheet1", now "Excel")
$ws.Name = "Excel"

# Make sure it's the active sheet, then move the selection to A2
$ws.Activate()
$ws.Range("A2").Select()
